$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Extend the repeating "Assignment / Grade / Comments" block with a new
#     7th block in columns AC:AF (AC is the blank divider column), mirroring
#     the existing blocks. Copy formatting from the last existing block
#     (Y:AB) so borders/fills/fonts match exactly (style indices are not
#     exposed directly, so we clone via copy/paste-formats).
$ws.Range("Y1:AB1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Y2:AB14").Copy()
$ws.Range("AC2").PasteSpecial(-4122)   # xlPasteFormats

# --- New block header row (row 1): Assignment / Grade / Comments
$ws.Range("AD1").Value = "Assignment"
$ws.Range("AE1").Value = "Grade"
$ws.Range("AF1").Value = "Comments"

# --- Fill in the new grades for Gal Erez (row 2)
# 4 UE: not submitted yet
$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = "Not submitted"

# 5 UE: not submitted yet
$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = "Not submitted"

# 6 UE: not submitted yet
$ws.Range("Z2").Value = "6 UE"
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = "Not submitted"

# 7 CPP: graded
$ws.Range("AD2").Value = "7 CPP"
$ws.Range("AE2").Value = 92
$ws.Range("AF2").Value = "Very good!"

# --- Header rows wrap onto two lines now, so they need extra height
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 30

# --- Restore view state: scrolled right to show the newest columns, with
#     the last active cell/selection parked at T25
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("T25").Select()

Write-Host "Edit applied"
